$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments: new col A, resized col B and C ---
$ws.Columns.Item(1).ColumnWidth = 5.5
$ws.Columns.Item(2).ColumnWidth = 27.5
$ws.Columns.Item(3).ColumnWidth = 6.5

# --- Updated sensitivity ratio values (D2:AG4) ---
$ws.Range("D2").Value = 13338.372247000423
$ws.Range("E2").Value = 13283.944155104331
$ws.Range("F2").Value = 13238.734373095178
$ws.Range("G2").Value = 13210.951907617329
$ws.Range("H2").Value = 12968.092820979931
$ws.Range("I2").Value = 12840.463832470257
$ws.Range("J2").Value = 12685.605914831272
$ws.Range("K2").Value = 12510.299469997506
$ws.Range("L2").Value = 12389.028237889952
$ws.Range("M2").Value = 12209.541411352271
$ws.Range("N2").Value = 12096.415385677963
$ws.Range("O2").Value = 11976.947297568846
$ws.Range("P2").Value = 11835.492085113468
$ws.Range("Q2").Value = 11669.342760262096
$ws.Range("R2").Value = 11475.477260362937
$ws.Range("S2").Value = 11250.512872213571
$ws.Range("T2").Value = 10986.572324960682
$ws.Range("U2").Value = 10682.324471817494
$ws.Range("V2").Value = 10333.341026754411
$ws.Range("W2").Value = 9934.5980465126759
$ws.Range("X2").Value = 9480.4070969217155
$ws.Range("Y2").Value = 8964.3403799899788
$ws.Range("Z2").Value = 8379.1500837048352
$ws.Range("AA2").Value = 7716.6826377113293
$ws.Range("AB2").Value = 6967.7891841005403
$ws.Range("AC2").Value = 6122.2344902478308
$ws.Range("AD2").Value = 5168.6078627588404
$ws.Range("AE2").Value = 4094.2415416793278
$ws.Range("AF2").Value = 2885.1448076530419
$ws.Range("AG2").Value = 1525.9659688175734
$ws.Range("D3").Value = 13338.542825044666
$ws.Range("E3").Value = 13284.219713187093
$ws.Range("F3").Value = 13239.204696775707
$ws.Range("G3").Value = 13211.731646148211
$ws.Range("H3").Value = 12969.309963252561
$ws.Range("I3").Value = 12842.271250366895
$ws.Range("J3").Value = 12688.176019090655
$ws.Range("K3").Value = 12513.823277396727
$ws.Range("L3").Value = 12393.714438426317
$ws.Range("M3").Value = 12215.615435425159
$ws.Range("N3").Value = 12063.435671030664
$ws.Range("O3").Value = 11901.547002629593
$ws.Range("P3").Value = 11714.30108598122
$ws.Range("Q3").Value = 11499.066203670303
$ws.Range("R3").Value = 11252.991051842218
$ws.Range("S3").Value = 10972.988789691737
$ws.Range("T3").Value = 10651.64036392536
$ws.Range("U3").Value = 10288.281868093256
$ws.Range("V3").Value = 9879.4196755093744
$ws.Range("W3").Value = 9421.304841206671
$ws.Range("X3").Value = 8909.9551721896441
$ws.Range("Y3").Value = 8341.1935084501674
$ws.Range("Z3").Value = 7710.707733071984
$ws.Range("AA3").Value = 7014.1396678735382
$ws.Range("AB3").Value = 6247.2121316767334
$ws.Range("AC3").Value = 5405.9061909128086
$ws.Range("AD3").Value = 4486.7042162969701
$ws.Range("AE3").Value = 3486.9190420525938
$ws.Range("AF3").Value = 2405.1356657484262
$ws.Range("AG3").Value = 1241.8000145501048
$ws.Range("D4").Value = 13338.949252046166
$ws.Range("E4").Value = 13284.876270400009
$ws.Range("F4").Value = 13240.325311354542
$ws.Range("G4").Value = 13213.589486613329
$ws.Range("H4").Value = 12972.209981754071
$ws.Range("I4").Value = 12846.577686399947
$ws.Range("J4").Value = 12694.29966631796
$ws.Range("K4").Value = 12522.219261073662
$ws.Range("L4").Value = 12404.879992826036
$ws.Range("M4").Value = 12230.087681263713
$ws.Range("N4").Value = 12048.749716538106
$ws.Range("O4").Value = 11853.770367190045
$ws.Range("P4").Value = 11629.234380051092
$ws.Range("Q4").Value = 11372.30171905596
$ws.Range("R4").Value = 11079.98184857552
$ws.Range("S4").Value = 10749.152588178242
$ws.Range("T4").Value = 10372.509380286841
$ws.Range("U4").Value = 9949.7128993002007
$ws.Range("V4").Value = 9477.8827683710169
$ws.Range("W4").Value = 8954.2739921741886
$ws.Range("X4").Value = 8376.4309885782186
$ws.Range("Y4").Value = 7742.39597517947
$ws.Range("Z4").Value = 7050.9876239852229
$ws.Range("AA4").Value = 6302.1702848338537
$ws.Range("AB4").Value = 5497.5396612356381
$ws.Range("AC4").Value = 4640.9579331356308
$ws.Range("AD4").Value = 3739.380387659211
$ws.Range("AE4").Value = 2803.9271948548862
$ws.Range("AF4").Value = 1851.2687693788191
$ws.Range("AG4").Value = 905.41212695017452

# Re-apply the original number format/style (quote-prefixed integer format)
# to the cells we just wrote, since assigning .Value resets quote-prefix state.
$ws.Range("AH2").Copy()
$ws.Range("D2:AG4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
